# Edit script for NEW_HAMPSHIRE_2015.xlsx
# - Rename header row to machine-friendly column names
# - Title-case "de" -> "De" in a set of specific place names
# - Delete trailing metadata/footer rows (96-100)
# - Dimension will shrink automatically after row deletion

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case fixes: "de" -> "De" within specific cells ---
$ws.Range("B12").Value = "Villa De Álvarez"
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("B29").Value = "Acapulco De Juárez"
$ws.Range("B31").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B35").Value = "Mártir De Cuilapan"
$ws.Range("B38").Value = "Tlapa De Comonfort"
$ws.Range("B41").Value = "Pachuca De Soto"
$ws.Range("B42").Value = "Progreso De Obregón"
$ws.Range("B46").Value = "Tepatitlán De Morelos"
$ws.Range("B49").Value = "Unión De Tula"
$ws.Range("B60").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B65").Value = "Tlacolula De Matamoros"
$ws.Range("B69").Value = "Cadereyta De Montes"
$ws.Range("B71").Value = "Landa De Matamoros"
$ws.Range("B91").Value = "Noria De Ángeles"

# --- Delete trailing footer/metadata rows (96-100) ---
$ws.Range("A96:A100").EntireRow.Delete()

# --- Update dimension to reflect new used range ---
$ws.UsedRange | Out-Null
